$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Price" (D) and "Volume(1h)" (E) columns with refreshed crypto data.
$ws.Range('D2').Value = "'20.281.48"
$ws.Range('E2').Value = '  +1.36%  '
$ws.Range('D3').Value = "'1.451.85"
$ws.Range('E3').Value = '  +2.83%  '
$ws.Range('D5').Value = "'0.9392"
$ws.Range('E5').Value = '  -6.20%  '
$ws.Range('D6').Value = "'273.38"
$ws.Range('E6').Value = '  -0.70%  '
$ws.Range('D7').Value = "'0.3640"
$ws.Range('E7').Value = '  -0.77%  '
$ws.Range('D8').Value = "'0.3052"
$ws.Range('E8').Value = '  -2.06%  '
$ws.Range('D9').Value = "'39.74"
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E10').Value = '  -0.60%  '
$ws.Range('D11').Value = "'0.06526"
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('E12').Value = '  -0.20%  '
$ws.Range('D13').Value = "'5.359"
$ws.Range('E13').Value = '  -2.34%  '
$ws.Range('E14').Value = '  -0.43%  '
$ws.Range('D15').Value = "'6.078"
$ws.Range('E15').Value = '  -1.83%  '
$ws.Range('D16').Value = "'0.00001017"
$ws.Range('E16').Value = '  -0.44%  '
$ws.Range('D17').Value = "'1.446.88"
$ws.Range('E17').Value = '  +2.46%  '
$ws.Range('D18').Value = "'0.9564"
$ws.Range('E18').Value = '  -4.43%  '
$ws.Range('D19').Value = "'0.05711"
$ws.Range('E19').Value = '  +0.40%  '
$ws.Range('D20').Value = "'69.07"
$ws.Range('E20').Value = '  -2.65%  '
$ws.Range('D21').Value = "'5.399"
$ws.Range('E21').Value = '  -4.13%  '
$ws.Range('D22').Value = "'14.32"
$ws.Range('E22').Value = '  -2.67%  '
$ws.Range('D23').Value = "'10.81"
$ws.Range('E23').Value = '  -2.38%  '
$ws.Range('D24').Value = "'2.233"
$ws.Range('E24').Value = '  -1.25%  '
$ws.Range('D25').Value = "'20.321.88"
$ws.Range('E25').Value = '  +1.53%  '
$ws.Range('D26').Value = "'140.82"
$ws.Range('D27').Value = "'2.087"
$ws.Range('E27').Value = '  -8.13%  '
$ws.Range('E28').Value = '  -0.64%  '
$ws.Range('D29').Value = "'1.599.71"
$ws.Range('E29').Value = '  +1.85%  '
$ws.Range('D30').Value = "'111.15"
$ws.Range('D31').Value = "'3.940"
$ws.Range('E31').Value = '  -0.52%  '
$ws.Range('D32').Value = "'4.818"
$ws.Range('E32').Value = '  -9.18%  '
$ws.Range('D33').Value = "'0.7844"
$ws.Range('E33').Value = '  -4.95%  '
$ws.Range('D34').Value = "'0.07736"
$ws.Range('E34').Value = '  +0.72%  '
$ws.Range('D35').Value = "'1.492"
$ws.Range('E35').Value = '  +0.51%  '
$ws.Range('D36').Value = "'0.05632"
$ws.Range('E36').Value = '  -4.83%  '
$ws.Range('D37').Value = "'4.646"
$ws.Range('E37').Value = '  -5.35%  '
$ws.Range('D38').Value = "'1.121"
$ws.Range('E38').Value = '  +2.21%  '
$ws.Range('D39').Value = "'0.02006"
$ws.Range('E39').Value = '  -3.43%  '
$ws.Range('D40').Value = "'0.9465"
$ws.Range('E40').Value = '  -5.37%  '
$ws.Range('D41').Value = "'10.20"
$ws.Range('E41').Value = '  -3.06%  '
$ws.Range('D42').Value = "'0.1852"
$ws.Range('E42').Value = '  -2.88%  '
$ws.Range('D43').Value = "'7.340"
$ws.Range('E43').Value = '  -12.85%  '
$ws.Range('D44').Value = "'0.5232"
$ws.Range('E44').Value = '  -1.61%  '
$ws.Range('D45').Value = "'3.471"
$ws.Range('E45').Value = '  -1.76%  '
$ws.Range('D46').Value = "'11.76"
$ws.Range('E46').Value = '  -4.40%  '
$ws.Range('D47').Value = "'116.67"
$ws.Range('E47').Value = '  +0.26%  '
$ws.Range('D48').Value = "'0.5104"
$ws.Range('E48').Value = '  -1.99%  '
$ws.Range('D49').Value = "'1.736"
$ws.Range('E49').Value = '  -1.95%  '
$ws.Range('D50').Value = "'0.06386"
$ws.Range('E50').Value = '  +3.06%  '
$ws.Range('D51').Value = "'0.9855"
$ws.Range('E51').Value = '  -1.46%  '

# Reset style on the Price cells we just touched so no extra quote-prefix /
# number-format style sticks around on them (keeps formatting identical to before).
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
